$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Old values (for reference): rsidR = 5DE4E2FEB51F49D59B7F44927DE4D21E,
# bookmark id = 109154115907730560796741218032147284857
$newRsid = '57A75031DCEE472CB120400FAFFD293E'
$newBmId = '41378587828078598767060583618226580861'

# Paragraph 1 : "Test link before bookmark : <field REF bookmark1>"
$p1 = $d.Paragraphs(2)
$r1 = $p1.Range
$xml1 = '<w:p ' + $wNs + ' w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidP="009168BC">' +
          '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr>' +
          '<w:r><w:t xml:space="preserve">Test link before bookmark : </w:t></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="begin"/></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="separate"/></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p>'
$null = $r1.InsertXML($xml1)

# Paragraph 2 : "Test bookmark : bookmarked content" (bookmarkStart/End id change)
$p2 = $d.Paragraphs(3)
$r2 = $p2.Range
$xml2 = '<w:p ' + $wNs + ' w:rsidR="00C52979" w:rsidRDefault="00E02A2B" w:rsidP="009168BC">' +
          '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr>' +
          '<w:r><w:t>Test</w:t></w:r>' +
          '<w:r w:rsidR="00C52979"><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:r><w:t>bookmark</w:t></w:r>' +
          '<w:r w:rsidR="00C52979"><w:t xml:space="preserve"> : </w:t></w:r>' +
          '<w:bookmarkStart w:name="bookmark1" w:id="' + $newBmId + '"/>' +
          '<w:r><w:t>bookmarked content</w:t></w:r>' +
          '<w:bookmarkEnd w:id="' + $newBmId + '"/>' +
        '</w:p>'
$null = $r2.InsertXML($xml2)

# Paragraph 3 : "Test link after bookmark : <field REF bookmark1> "
$p3 = $d.Paragraphs(4)
$r3 = $p3.Range
$xml3 = '<w:p ' + $wNs + ' w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidP="00E02A2B">' +
          '<w:pPr><w:tabs><w:tab w:val="left" w:pos="3119"/></w:tabs></w:pPr>' +
          '<w:r><w:t xml:space="preserve">Test link after bookmark : </w:t></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="begin"/></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="separate"/></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +
          '<w:r w:rsidR="' + $newRsid + '"><w:fldChar w:fldCharType="end"/></w:r>' +
          '<w:r w:rsidR="00D0546C"><w:t xml:space="preserve"> </w:t></w:r>' +
        '</w:p>'
$null = $r3.InsertXML($xml3)

Write-Host "Done"
